$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Area - 146.01*116.81"
$ws.Range("F3").Formula = "=146.01*116.81"

$ws.Range("E4").Value = "STDEV"
$ws.Range("F4").Formula = "=STDEV(B2:B13)"

$ws.Range("E5").Value = "Mean"
$ws.Range("F5").Value = 397.17

$ws.Range("E6").Value = "Area Fraction"
$ws.Range("F6").Formula = "=(SUM(B2:B13))*100/F3"

$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

$ws.Range("B2").Select()
